$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.91194
$ws.Range("B3").Value = 1.91194
$ws.Range("B4").Value = 1.91194
$ws.Range("B5").Value = 1.91194
$ws.Range("B6").Value = 1.5278
$ws.Range("B7").Value = 1.5278
$ws.Range("B8").Value = 1.91194
$ws.Range("B9").Value = 1.91194
$ws.Range("B10").Value = 1.91194
$ws.Range("B11").Value = 1.91194
$ws.Range("B12").Value = 1.57672
$ws.Range("B13").Value = 1.5278
$ws.Range("B14").Value = 1.91194
$ws.Range("B15").Value = 1.91194
$ws.Range("B16").Value = 1.91194
$ws.Range("B17").Value = 1.57672
$ws.Range("B18").Value = 1.5278
$ws.Range("B19").Value = 1.91194
$ws.Range("B20").Value = 1.91194
$ws.Range("B21").Value = 1.91194
$ws.Range("B22").Value = 1.5278
$ws.Range("B23").Value = 1.5278
$ws.Range("B24").Value = 1.5278
$ws.Range("B25").Value = 1.91194
$ws.Range("B26").Value = 1.91194
$ws.Range("B27").Value = 1.5278
$ws.Range("B28").Value = 1.5278
$ws.Range("B29").Value = 1.5278
$ws.Range("B30").Value = 1.91194
$ws.Range("B31").Value = 1.91194
$ws.Range("B32").Value = 1.5278
$ws.Range("B33").Value = 1.5278
$ws.Range("B34").Value = 1.5278
$ws.Range("B35").Value = 1.5278
$ws.Range("B36").Value = 1.91194
$ws.Range("B37").Value = 1.91194
$ws.Range("B38").Value = 1.91194
$ws.Range("B39").Value = 1.88589
$ws.Range("B40").Value = 1.5278
$ws.Range("B41").Value = 1.5278
$ws.Range("B42").Value = 1.91194
$ws.Range("B43").Value = 1.91194
$ws.Range("B44").Value = 1.5278
$ws.Range("B45").Value = 1.5278
$ws.Range("B46").Value = 1.5278
$ws.Range("B47").Value = 1.91194
$ws.Range("B48").Value = 1.91194
$ws.Range("B49").Value = 1.91194
$ws.Range("B50").Value = 1.5278
$ws.Range("B51").Value = 1.5278
$ws.Range("B52").Value = 1.5278
$ws.Range("B53").Value = 1.91194
$ws.Range("B54").Value = 1.91194
$ws.Range("B55").Value = 1.64399
$ws.Range("B56").Value = 1.5278
$ws.Range("B57").Value = 1.5278
$ws.Range("B58").Value = 1.5278
$ws.Range("B59").Value = 1.91194
$ws.Range("B60").Value = 1.5278
$ws.Range("B61").Value = 1.5278
$ws.Range("B62").Value = 1.5278
$ws.Range("B63").Value = 1.5278
$ws.Range("B64").Value = 1.91194
$ws.Range("B65").Value = 1.91194
$ws.Range("B66").Value = 1.5278
$ws.Range("B67").Value = 1.5278
$ws.Range("B68").Value = 1.5278
